$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text number format only for Price (D) column cells so Excel keeps
# the exact text representation (e.g. thousands-dot notation, leading zeros)
# instead of re-interpreting the value as a number.
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D50').NumberFormat = '@'

$ws.Range('D2').Value = '57.693.44'
$ws.Range('E2').Value = '  -1.23%  '
$ws.Range('D3').Value = '2.564.56'
$ws.Range('E3').Value = '  -3.14%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').Value = '515.26'
$ws.Range('E5').Value = '  -1.38%  '
$ws.Range('D6').Value = '141.95'
$ws.Range('E6').Value = '  -2.04%  '
$ws.Range('E7').Value = '  -0.15%  '
$ws.Range('D8').Value = '0.563'
$ws.Range('E8').Value = '  -1.63%  '
$ws.Range('D9').Value = '2.579.33'
$ws.Range('E9').Value = '  -2.77%  '
$ws.Range('D10').Value = '6.61'
$ws.Range('E10').Value = '  -1.88%  '
$ws.Range('E11').Value = '  -2.28%  '
$ws.Range('E12').Value = '  -4.93%  '
$ws.Range('E13').Value = '  -1.14%  '
$ws.Range('D14').Value = '3.016.01'
$ws.Range('E14').Value = '  -3.17%  '
$ws.Range('D15').Value = '57.716.38'
$ws.Range('E15').Value = '  -1.20%  '
$ws.Range('D16').Value = '20.18'
$ws.Range('E16').Value = '  -3.65%  '
$ws.Range('B17').Value = 'WrappedEther'
$ws.Range('C17').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D17').Value = '2.596.43'
$ws.Range('E17').Value = '  -2.16%  '
$ws.Range('B18').Value = 'ShibaInu'
$ws.Range('C18').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D18').Value = '0.0000133'
$ws.Range('E18').Value = '  -2.63%  '
$ws.Range('D19').Value = '337.26'
$ws.Range('E19').Value = '  -0.36%  '
$ws.Range('E20').Value = '  -2.65%  '
$ws.Range('D21').Value = '10.18'
$ws.Range('E21').Value = '  -2.73%  '
$ws.Range('E22').Value = '  -0.19%  '
$ws.Range('E23').Value = '  -0.36%  '
$ws.Range('D24').Value = '65.20'
$ws.Range('E24').Value = '  +1.36%  '
$ws.Range('E25').Value = '  -1.01%  '
$ws.Range('D26').Value = '0.996'
$ws.Range('E26').Value = '  -0.32%  '
$ws.Range('D27').Value = '0.399'
$ws.Range('E27').Value = '  -6.10%  '
$ws.Range('D28').Value = '2.683.54'
$ws.Range('E28').Value = '  -3.07%  '
$ws.Range('D29').Value = '6.95'
$ws.Range('E29').Value = '  -2.67%  '
$ws.Range('D30').Value = '0.0₃0740'
$ws.Range('E30').Value = '  -7.31%  '
$ws.Range('E31').Value = '  -0.06%  '
$ws.Range('D32').Value = '6.22'
$ws.Range('E32').Value = '  -7.19%  '
$ws.Range('E33').Value = '  -1.41%  '
$ws.Range('B34').Value = 'EthereumClassic'
$ws.Range('C34').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D34').Value = '18.65'
$ws.Range('E34').Value = '  -1.29%  '
$ws.Range('B35').Value = 'Monero'
$ws.Range('C35').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D35').Value = '149.88'
$ws.Range('E35').Value = '  -1.89%  '
$ws.Range('E36').Value = '  -4.04%  '
$ws.Range('E37').Value = '  -4.22%  '
$ws.Range('D38').Value = '0.863'
$ws.Range('E38').Value = '  -5.05%  '
$ws.Range('D39').Value = '35.99'
$ws.Range('E39').Value = '  -2.44%  '
$ws.Range('E40').Value = '  -4.36%  '
$ws.Range('E41').Value = '  -0.90%  '
$ws.Range('D42').Value = '3.51'
$ws.Range('E42').Value = '  -3.29%  '
$ws.Range('D43').Value = '0.998'
$ws.Range('E43').Value = '  -0.24%  '
$ws.Range('D44').Value = '269.07'
$ws.Range('E44').Value = '  -1.55%  '
$ws.Range('D45').Value = '10.66'
$ws.Range('E45').Value = '  +0.32%  '
$ws.Range('D46').Value = '0.0948'
$ws.Range('E46').Value = '  -2.46%  '
$ws.Range('D47').Value = '0.584'
$ws.Range('E47').Value = '  -3.92%  '
$ws.Range('D48').Value = '18.68'
$ws.Range('E48').Value = '  -3.99%  '
$ws.Range('E49').Value = '  -2.94%  '
$ws.Range('D50').Value = '1.970.03'
$ws.Range('E50').Value = '  -3.72%  '
$ws.Range('E51').Value = '  -4.14%  '
